$wb = $excel.ActiveWorkbook

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 900.4286
$ws.Range("I32").Value = 533.6667
$ws.Range("J32").Value = 1175.5
$ws.Range("K32").Value = 533.6667
$ws.Range("L32").Value = 1175.5
$ws.Range("M32").Value = -207.6667
$ws.Range("N32").Value = -1827.5

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 250.59459
$ws.Range("I33").Value = 239.48276
$ws.Range("J33").Value = 290.875
$ws.Range("K33").Value = 239.48276
$ws.Range("L33").Value = 290.875
$ws.Range("M33").Value = -10.48276000000001
$ws.Range("N33").Value = -748.875

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3907.4075
$ws.Range("I74").Value = 3843.7144
$ws.Range("J74").Value = 3929.7
$ws.Range("K74").Value = 3843.7144
$ws.Range("L74").Value = 3929.7
$ws.Range("M74").Value = -2907.7144
$ws.Range("N74").Value = -5801.7

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3470.6
$ws.Range("I76").Value = 3375
$ws.Range("J76").Value = 3853
$ws.Range("K76").Value = 3375
$ws.Range("L76").Value = 3853
$ws.Range("M76").Value = -3060
$ws.Range("N76").Value = -4483

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3907.4075
$ws.Range("I77").Value = 3843.7144
$ws.Range("J77").Value = 3929.7
$ws.Range("K77").Value = 19218.572
$ws.Range("L77").Value = 19648.5
$ws.Range("M77").Value = -14538.572
$ws.Range("N77").Value = -29008.5

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3470.6
$ws.Range("I79").Value = 3375
$ws.Range("J79").Value = 3853
$ws.Range("K79").Value = 3375
$ws.Range("L79").Value = 3853
$ws.Range("M79").Value = -2283
$ws.Range("N79").Value = -6037

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4925.6665
$ws.Range("I113").Value = 4861
$ws.Range("J113").Value = 4971.857
$ws.Range("K113").Value = 4861
$ws.Range("L113").Value = 4971.857
$ws.Range("M113").Value = -1607
$ws.Range("N113").Value = -11479.857

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2600725.2
$ws.Range("I138").Value = 2619.077
$ws.Range("J138").Value = 3925250.2
$ws.Range("K138").Value = 7857.231000000001
$ws.Range("L138").Value = 11775750.6
$ws.Range("M138").Value = -2717.231000000001
$ws.Range("N138").Value = -11786030.6

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2706.652
$ws.Range("I2").Value = 2446.7144
$ws.Range("J2").Value = 3111
$ws.Range("K2").Value = 2446.7144
$ws.Range("L2").Value = 3111
$ws.Range("M2").Value = -2333.7144
$ws.Range("N2").Value = -3337

# ARM row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("N23").Value = 0

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4099.6
$ws.Range("I63").Value = 4121.0713
$ws.Range("J63").Value = 3799
$ws.Range("K63").Value = 4121.0713
$ws.Range("L63").Value = 3799
$ws.Range("M63").Value = -3435.0713
$ws.Range("N63").Value = -5171

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4099.6
$ws.Range("I66").Value = 4121.0713
$ws.Range("J66").Value = 3799
$ws.Range("K66").Value = 20605.3565
$ws.Range("L66").Value = 18995
$ws.Range("M66").Value = -17173.3565
$ws.Range("N66").Value = -25859

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 28573250
$ws.Range("I102").Value = 28573250
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 28573250
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -28571628

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2706.652
$ws.Range("I116").Value = 2446.7144
$ws.Range("J116").Value = 3111
$ws.Range("K116").Value = 2446.7144
$ws.Range("L116").Value = 3111
$ws.Range("M116").Value = -152.7143999999998
$ws.Range("N116").Value = -7699

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2706.652
$ws.Range("I3").Value = 2446.7144
$ws.Range("J3").Value = 3111
$ws.Range("K3").Value = 2446.7144
$ws.Range("L3").Value = 3111
$ws.Range("M3").Value = -2332.7144
$ws.Range("N3").Value = -3339

# BSM row 11
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 2069
$ws.Range("I11").Value = 938.6667
$ws.Range("J11").Value = 4329.6665
$ws.Range("K11").Value = 938.6667
$ws.Range("L11").Value = 4329.6665
$ws.Range("M11").Value = -798.6667
$ws.Range("N11").Value = -4609.6665

# CRP row 2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 33535
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 100005
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 100005
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -100231

# CUL row 101
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 8000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 8000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 24000
$ws.Range("N101").Value = -28868

# GSM row 4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 5916.6665
$ws.Range("I4").Value = 1750
$ws.Range("J4").Value = 8000
$ws.Range("K4").Value = 1750
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = -1638
$ws.Range("N4").Value = -8224

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2506106
$ws.Range("I5").Value = 5000504
$ws.Range("J5").Value = 11708
$ws.Range("K5").Value = 5000504
$ws.Range("L5").Value = 11708
$ws.Range("M5").Value = -5000392
$ws.Range("N5").Value = -11932

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3783.795
$ws.Range("I80").Value = 4020
$ws.Range("J80").Value = 3749.0588
$ws.Range("K80").Value = 4020
$ws.Range("L80").Value = 3749.0588
$ws.Range("M80").Value = -3022
$ws.Range("N80").Value = -5745.0588

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3783.795
$ws.Range("I83").Value = 4020
$ws.Range("J83").Value = 3749.0588
$ws.Range("K83").Value = 20100
$ws.Range("L83").Value = 18745.294
$ws.Range("M83").Value = -15108
$ws.Range("N83").Value = -28729.294

# GSM row 103
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2087
$ws.Range("I113").Value = 1400.4615
$ws.Range("J113").Value = 2898.3635
$ws.Range("K113").Value = 1400.4615
$ws.Range("L113").Value = 2898.3635
$ws.Range("M113").Value = 769.5385000000001
$ws.Range("N113").Value = -7238.363499999999

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 238.6923
$ws.Range("I55").Value = 232.71428
$ws.Range("J55").Value = 263.8
$ws.Range("K55").Value = 232.71428
$ws.Range("L55").Value = 263.8
$ws.Range("M55").Value = -59.71428
$ws.Range("N55").Value = -609.8

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1775.7059
$ws.Range("I61").Value = 1759.3043
$ws.Range("J61").Value = 1810
$ws.Range("K61").Value = 1759.3043
$ws.Range("L61").Value = 1810
$ws.Range("M61").Value = -1557.3043
$ws.Range("N61").Value = -2214

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1775.7059
$ws.Range("I113").Value = 1759.3043
$ws.Range("J113").Value = 1810
$ws.Range("K113").Value = 1759.3043
$ws.Range("L113").Value = 1810
$ws.Range("M113").Value = 410.6957
$ws.Range("N113").Value = -6150

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3848.3635
$ws.Range("I122").Value = 3803.7646
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 11411.2938
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8961.293799999999
$ws.Range("N122").Value = -16900

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1066666.6
$ws.Range("I4").Value = 1066666.6
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1066666.6
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -1066553.6

# WVR row 6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 558.1667
$ws.Range("I6").Value = 75
$ws.Range("J6").Value = 799.75
$ws.Range("K6").Value = 75
$ws.Range("L6").Value = 799.75
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = -1029.75

